# Updated cryptos list on Sat Sep  7 19:49:34 UTC 2024 with GitHub Actions
#
# Applies the latest scrape of coinranking.com data to the "Sheet1" table:
# most rows just get refreshed Price (col D) / Volume(1h) (col E) text, a
# couple of coins (rows 39/40) swapped rank order with new figures.
#
# Price/volume cells are stored as *text* (t="inlineStr" in the original
# OOXML) even when they look like plain numbers (e.g. "0.999", "495.86").
# Assigning such a look-alike string straight to Range.Value lets Excel's
# usual cell-content inference silently coerce it into a real number, so
# we force those through as text (leading "'" = quote-prefix) and then
# restore the cell to the default "Normal" style so no stray number
# format/style sticks around afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Value
    )

    $range = $ws.Range($Cell)

    if ($Value -match '^-?\d+(\.\d+)?$') {
        # Would otherwise be auto-converted to a number by Excel.
        $range.Value = "'" + $Value
        $range.Style = "Normal"
    } else {
        $range.Value = $Value
    }
}

Set-TextValue "D2"  "54.101.95"
Set-TextValue "E2"  "  +0.95%  "

Set-TextValue "D3"  "2.287.90"

Set-TextValue "E4"  "  +0.32%  "

Set-TextValue "D5"  "495.86"
Set-TextValue "E5"  "  +2.44%  "

Set-TextValue "D6"  "128.57"
Set-TextValue "E6"  "  +2.67%  "

Set-TextValue "D7"  "0.999"
Set-TextValue "E7"  "  -0.09%  "

Set-TextValue "E8"  "  +2.80%  "

Set-TextValue "D9"  "2.286.34"
Set-TextValue "E9"  "  +3.42%  "

Set-TextValue "D10" "0.0950"
Set-TextValue "E10" "  +4.59%  "

Set-TextValue "D11" "0.152"
Set-TextValue "E11" "  +2.52%  "

Set-TextValue "E12" "  +4.64%  "

Set-TextValue "E13" "  +0.13%  "

Set-TextValue "D14" "2.692.26"
Set-TextValue "E14" "  +3.18%  "

Set-TextValue "D15" "21.81"
Set-TextValue "E15" "  +4.42%  "

Set-TextValue "D16" "54.192.90"
Set-TextValue "E16" "  +1.26%  "

Set-TextValue "E17" "  +1.97%  "

Set-TextValue "D18" "2.285.44"
Set-TextValue "E18" "  +3.28%  "

Set-TextValue "D19" "10.05"
Set-TextValue "E19" "  +5.75%  "

Set-TextValue "E20" "  +4.72%  "

Set-TextValue "D21" "6.46"
Set-TextValue "E21" "  +6.25%  "

Set-TextValue "D22" "301.11"
Set-TextValue "E22" "  +1.37%  "

Set-TextValue "E23" "  -0.14%  "

Set-TextValue "E24" "  -1.90%  "

Set-TextValue "D25" "62.70"
Set-TextValue "E25" "  -0.64%  "

Set-TextValue "E26" "  +0.93%  "

Set-TextValue "E27" "  +3.36%  "

Set-TextValue "D28" "2.389.27"
Set-TextValue "E28" "  +2.72%  "

Set-TextValue "E29" "  +4.52%  "

Set-TextValue "D30" "7.07"
Set-TextValue "E30" "  +1.98%  "

Set-TextValue "D31" "169.10"
Set-TextValue "E31" "  -0.03%  "

Set-TextValue "D32" "1.61"
Set-TextValue "E32" "  +2.36%  "

Set-TextValue "D33" "0.0₃0690"
Set-TextValue "E33" "  +2.84%  "

Set-TextValue "D34" "5.89"
Set-TextValue "E34" "  +3.46%  "

Set-TextValue "E35" "  +0.04%  "

Set-TextValue "D36" "0.996"
Set-TextValue "E36" "  +0.33%  "

Set-TextValue "E37" "  +2.60%  "

Set-TextValue "D38" "17.69"
Set-TextValue "E38" "  +2.17%  "

# Rows 39/40 swapped order: SuiNetwork now ranks above ImmutableX.
Set-TextValue "B39" "SuiNetwork"
Set-TextValue "C39" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D39" "0.906"
Set-TextValue "E39" "  +9.76%  "

Set-TextValue "B40" "ImmutableX"
Set-TextValue "C40" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D40" "1.19"
Set-TextValue "E40" "  +4.71%  "

Set-TextValue "E41" "  +4.69%  "

Set-TextValue "D42" "35.49"
Set-TextValue "E42" "  -0.53%  "

Set-TextValue "E43" "  +3.79%  "

Set-TextValue "E44" "  +3.17%  "

Set-TextValue "E45" "  +3.65%  "

Set-TextValue "D46" "127.28"
Set-TextValue "E46" "  +3.99%  "

Set-TextValue "D47" "4.79"
Set-TextValue "E47" "  +3.59%  "

Set-TextValue "D48" "0.0888"
Set-TextValue "E48" "  +1.46%  "

Set-TextValue "D49" "0.546"
Set-TextValue "E49" "  +4.04%  "

Set-TextValue "D50" "238.75"
Set-TextValue "E50" "  +4.50%  "

Set-TextValue "E51" "  +3.84%  "
